$p = $ppt.ActivePresentation

# The deck was regenerated by the latest binary. The visible/meaningful
# effect of that regeneration is that placeholder shapes on the slide
# master, the slide layout and slide 1 now carry an explicit rectangle
# preset geometry (<a:prstGeom prst="rect"><a:avLst/></a:prstGeom>)
# instead of relying on inherited/implicit geometry. Setting
# AutoShapeType to msoShapeRectangle (1) on a shape materializes that
# geometry in the OOXML without otherwise touching the shape.

$msoShapeRectangle = 1

# --- Slide master placeholders (Title, Content, Page Numbering) ---
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shape = $master.Shapes.Item($i)
    if ($shape.Name -eq "Title" -or $shape.Name -eq "Content" -or $shape.Name -eq "Page Numbering") {
        $shape.AutoShapeType = $msoShapeRectangle
    }
}

# --- Slide layout placeholders (Title, Content, Page Numbering) ---
$layout = $master.CustomLayouts.Item(1)
for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
    $shape = $layout.Shapes.Item($i)
    if ($shape.Name -eq "Title" -or $shape.Name -eq "Content" -or $shape.Name -eq "Page Numbering") {
        $shape.AutoShapeType = $msoShapeRectangle
    }
}

# --- Slide 1 placeholders (Title, Page Numbering) ---
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.Name -eq "Title" -or $shape.Name -eq "Page Numbering") {
        $shape.AutoShapeType = $msoShapeRectangle
    }
}
